$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rows 2-49: bump the rank number in column A by 1 (it was 0-based, now
#    1-based).
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 49; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 1
}

# ---------------------------------------------------------------------------
# 2. Row 51: clarify the totals label.
# ---------------------------------------------------------------------------
$ws.Range("B51").Value = "Total's for top 50 (actually 48)"

# ---------------------------------------------------------------------------
# 3. "Old Whales" block (rows 54-56): the wallets got re-sorted and two of
#    the balances were refreshed. New row order:
#      54: 0x2bd6997b...  (previously row 56)
#      55: 0x3d268cd5...  (previously row 54, refreshed amount)
#      56: 0x505dd22c...  (previously row 55, refreshed amount)
# ---------------------------------------------------------------------------

# -- Row 54 ----------------------------------------------------------------
$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "0x2bd6997bf6fcfde139eb1b9346fbf79defd4e8cc"
$ws.Range("C54").Value = "'0"
$ws.Range("D54").Value = "0 Trillion"
$ws.Range("E54").Value = "0.5497% "
$ws.Range("F54").Value = "550 Trillion Less"
$ws.Range("G54").Value = "'549,662,717,019,826"

# -- Row 55 ----------------------------------------------------------------
$ws.Range("A55").Value = 2
$ws.Range("B55").Value = "0x3d268cd580f89cfe6cc5dcf8764f51085f74a649"
$ws.Range("C55").Value = "'4,787,982,343,926"
$ws.Range("D55").Value = "5 Trillion"
$ws.Range("E55").Value = "0.1278% "
$ws.Range("F55").Value = "123 Trillion Less"
$ws.Range("G55").Value = "'123,043,485,437,911"

# -- Row 56 ----------------------------------------------------------------
$ws.Range("A56").Value = 3
$ws.Range("B56").Value = "0x505dd22c1bacced7531f319f5008318a440490bc"
$ws.Range("C56").Value = "'90,051,450,064,196"
$ws.Range("D56").Value = "90 Trillion"
$ws.Range("E56").Value = "0.0950% "
$ws.Range("F56").Value = "5 Trillion Less"
$ws.Range("G56").Value = "'4,952,796,909,727"

# Writing digit-only text through .Value (with the leading quote needed to
# keep it text) bumps the cell style to a quote-prefixed variant, so stamp
# the original "old whales" formatting back on using row 15 (same style set:
# s4/s4/s5/s5/s5/s6/s7) as a pristine, untouched format source.
$ws.Range("A15:G15").Copy()
$ws.Range("A54:G54").PasteSpecial(-4122)
$ws.Range("A15:G15").Copy()
$ws.Range("A55:G55").PasteSpecial(-4122)
$ws.Range("A15:G15").Copy()
$ws.Range("A56:G56").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. New row 58: grand total loss / gain across all whales.
# ---------------------------------------------------------------------------
$ws.Range("E58").Value = "Total loss / gain from all whales"
$ws.Range("F58").Value = "-689 Trillion"
$ws.Range("G58").Value = "'-688,641,404,522,513"

$ws.Range("F54").Copy()
$ws.Range("D58").PasteSpecial(-4122)
$ws.Range("G54").Copy()
$ws.Range("E58").PasteSpecial(-4122)
$ws.Range("F54").Copy()
$ws.Range("F58").PasteSpecial(-4122)
$ws.Range("F54").Copy()
$ws.Range("G58").PasteSpecial(-4122)
